# Add 2022-Q3 data
# -----------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 right after
#    the header, shifting the existing quarters down, and renumber the
#    index column (A) sequentially.
# -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")

$summary.Rows("2:2").Insert()
$summary.Rows("2:2").ClearFormats()

# restore the index-column style (bold/centered/bordered) that row 2
# should carry, by copying it from the row right below (still holding
# the old row-2 formatting)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

# renumber the remaining index column sequentially
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# -----------------------------------------------------------------
# 2) Create the new "2022-Q3" sheet (positioned right after "总计"),
#    using the existing "2021-Q4" sheet as a formatting template, then
#    replace its data with the new fund holdings.
# -----------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $summary)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# drop the extra template data rows, keep header (row 1) + 2 data rows
$newSheet.Rows("4:11").Delete()

# force the numeric-looking text fields to stay text (keep leading
# zeros on fund codes, and preserve the original text formatting of
# percentages/amounts), matching the other quarter sheets
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "090019"
$newSheet.Range("C2").Value = "大成景恒混合A"
$newSheet.Range("D2").Value = "1.13"
$newSheet.Range("E2").Value = "93.98"
$newSheet.Range("F2").Value = "1.79"
$newSheet.Range("G2").Value = "0.0202"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006038"
$newSheet.Range("C3").Value = "大成景恒混合C"
$newSheet.Range("D3").Value = "0.45"
$newSheet.Range("E3").Value = "93.98"
$newSheet.Range("F3").Value = "1.79"
$newSheet.Range("G3").Value = "0.0081"
$newSheet.Range("H3").Value = 7

# restore the originally-selected tab (last sheet) as the active one,
# same as before the edit
$wb.Worksheets.Item("2021-Q2").Select()
